$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.014.55"
$ws.Range("E2").Value = "  +1.83%  "
$ws.Range("D3").Value = "2.502.12"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'591.10"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").Value = "'174.88"
$ws.Range("E6").Value = "  -0.82%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.514"
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("D9").Value = "2.501.35"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("E10").Value = "  +5.82%  "
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("D12").Value = "'4.95"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("E13").Value = "  -1.83%  "
$ws.Range("D14").Value = "2.946.69"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").Value = "'25.55"
$ws.Range("E15").Value = "  -0.93%  "
$ws.Range("D16").Value = "68.882.46"
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "2.484.54"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("D19").Value = "'360.38"
$ws.Range("E19").Value = "  +2.62%  "
$ws.Range("E20").Value = "  -0.55%  "
$ws.Range("E21").Value = "  -2.26%  "
$ws.Range("E22").Value = "  -2.14%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "'69.99"
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("E26").Value = "  -3.35%  "
$ws.Range("E27").Value = "  -7.77%  "
$ws.Range("D28").Value = "2.630.21"
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("E29").Value = "  +0.29%  "
$ws.Range("D30").Value = "'503.07"
$ws.Range("E30").Value = "  -0.94%  "
$ws.Range("D31").Value = "0.0₃0875"
$ws.Range("E31").Value = "  -3.58%  "
$ws.Range("D32").Value = "'7.68"
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("E33").Value = "  -0.86%  "
$ws.Range("E34").Value = "  -4.78%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "'162.43"
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("E37").Value = "  -4.11%  "
$ws.Range("D38").Value = "'18.58"
$ws.Range("E38").Value = "  +0.97%  "
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").Value = "'1.29"
$ws.Range("E41").Value = "  -3.86%  "
$ws.Range("E42").Value = "  -3.07%  "
$ws.Range("D43").Value = "'4.71"
$ws.Range("E43").Value = "  -2.84%  "
$ws.Range("D44").Value = "'0.318"
$ws.Range("E44").Value = "  -3.76%  "
$ws.Range("E45").Value = "  -5.33%  "
$ws.Range("D46").Value = "'149.13"
$ws.Range("E46").Value = "  +2.53%  "
$ws.Range("E47").Value = "  +0.38%  "
$ws.Range("E48").Value = "  -1.59%  "
$ws.Range("E49").Value = "  -1.37%  "
$ws.Range("E50").Value = "  -2.49%  "
